$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The sheet is a "label / value" table (col A = label, col B/C = value).
# A new row needs to appear right after "Docentes responsaveis:" (row 12) so
# the professor's name (previously mis-placed beside "Objetivos:") becomes
# its own value row, with every row below it shifting down by one. Several
# cells also receive brand new text that was missing before.
# ---------------------------------------------------------------------------

# 1) Insert a new blank row at row 13 - this shifts old rows 13..23 down to
#    14..24, carrying their formatting (row height, styles, shared content).
$ws.Rows.Item(13).Insert()

# The Insert() leaves a stray formatted-but-empty cell at A13 (copied from
# A12's bold style). Clear it completely so row 13 has no A cell at all,
# matching the target layout.
$ws.Range("A13").Clear()

# 2) Row 10 ("Objetivos:") gets a brand new objectives paragraph - it used
#    to wrongly hold the professor's name.
$newObjetivos = "Complementar a formação multidisciplinar dos alunos de Engenharia abordando, com maior profundidade, tópicos atuais e relevantes sobre fenômenos de transporte, termodinâmica, operações unitárias e reatore"
$ws.Range("B10").Value = $newObjetivos
$ws.Range("C10").Value = $newObjetivos

# 3) Row 13 (the new row, no label) now correctly holds the professor's name
#    that belongs to "Docentes responsáveis:" (row 12).
$docente = "5816812 - João Paulo Alves Silva"
$ws.Range("B13").Value = $docente
$ws.Range("C13").Value = $docente
# Give the brand-new B13/C13 cells the same formatting as their column
# siblings (Insert() only guessed a style based on column A's format).
$ws.Range("B14").Copy()
$ws.Range("B13").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("C14").Copy()
$ws.Range("C13").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# 4) Row 14 ("Programa resumido:") gets the real short-syllabus text instead
#    of the placeholder "Semestral".
$programaResumido = "Tópicos atuais e relevantes sobre fenômenos de transporte, termodinâmica, operações unitárias e reatores."
$ws.Range("B14").Value = $programaResumido
$ws.Range("C14").Value = $programaResumido

# 5) Row 16 ("Programa:") gets the real syllabus text instead of the
#    duplicated activation date.
$programa = "O conteúdo desta disciplina será de acordo com o tópico a ser programado, devendo abordar assuntos complementares a formação de um profissional de Engenharia Química."
$ws.Range("B16").Value = $programa
$ws.Range("C16").Value = $programa

# 6) Row 19 ("Método:") gets the "development of the course" text (it was
#    wrongly holding the professor's name after the shift).
$metodo = "O desenvolvimento da disciplina será baseado em leituras, aula expositiva, discussão e resolução de estudos de caso e resolução de exercícios."
$ws.Range("B19").Value = $metodo
$ws.Range("C19").Value = $metodo

# 7) Row 20 ("Critério:") gets "Provas e trabalhos."
$criterio = "Provas e trabalhos."
$ws.Range("B20").Value = $criterio
$ws.Range("C20").Value = $criterio

# 8) Row 21 ("Norma de recuperação:") gets the recovery-exam rule text.
$norma = "Prova única com nota maior ou igual a 5,0 (cinco)."
$ws.Range("B21").Value = $norma
$ws.Range("C21").Value = $norma

# 9) Row 22 ("Bibliografia:") gets the real bibliography text (two lines).
$bibliografia = "Textos fornecidos pelo professor da disciplina" + [char]10 + "Artigos extraídos de revistas especializadas de Engenharia Química."
$ws.Range("B22").Value = $bibliografia
$ws.Range("C22").Value = $bibliografia

# Rows 11, 12, 15, 17, 18, 23, 24 already carry the correct content after the
# row-13 insert shifted everything down, so nothing else needs to change.
